$wb = $excel.ActiveWorkbook

# Both the "展览" sheet (sheet1) and the "全部类型" sheet (sheet4) contain
# identical copies of this data; update the "想去人数" (F column) values
# for the affected event rows on both sheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F4").Value = 63
    $ws.Range("F5").Value = 375
    $ws.Range("F6").Value = 11132
    $ws.Range("F7").Value = 578
    $ws.Range("F18").Value = 314
    $ws.Range("F19").Value = 1195
    $ws.Range("F20").Value = 59
}
